# The document ends with three list paragraphs:
#   "Hover over to show full name of article"   (ListParagraph, ilvl=1)
#   <empty, but holds the _GoBack bookmark>       (ListParagraph, ilvl=1)
#   <empty>                                       (ListParagraph, ilvl=0)
#
# The edit merges all three into a single paragraph: the "Hover over..."
# text immediately followed by the (still-empty) _GoBack bookmark, with no
# trailing empty list paragraphs left afterwards.
#
# We do this the same way a user would in Word: put the cursor at the end
# of the "Hover over..." paragraph and press Delete twice, merging the two
# following (empty) paragraphs into it one at a time. Deleting just the
# paragraph mark (rather than the whole paragraph range) keeps any
# bookmarks that live inside the paragraph being merged forward.

$d = $word.ActiveDocument

# Locate the "Hover over..." paragraph by its text so the script is
# resilient to the exact paragraph index.
$targetText = "Hover over to show full name of article"
$index = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext.TrimEnd([char]13) -eq $targetText) {
        $index = $i
    }
}

# First merge: delete the paragraph mark ending the "Hover over..."
# paragraph, pulling the bookmark-only paragraph after it into place.
$p = $d.Paragraphs.Item($index)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

# Second merge: delete the paragraph mark again to pull the final, wholly
# empty paragraph into place too, leaving a single merged paragraph.
$p = $d.Paragraphs.Item($index)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

# Each merge keeps the *later* paragraph's formatting (matching Word's own
# paragraph-mark-delete behavior), so the surviving paragraph mark came
# from the last (ilvl=0) paragraph. Restore the original "Hover over..."
# list level (ilvl=1, i.e. ListLevelNumber 2) that should remain on the
# merged paragraph.
$final = $d.Paragraphs.Item($index)
$final.Range.ListFormat.ListLevelNumber = 2
